$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 550
$ws.Range("H21").Value = 225008.5
$ws.Range("I21").Value = 225008.5
$ws.Range("K21").Value = 225008.5
$ws.Range("M21").Value = -224540.5
$ws.Range("H23").Value = 225008.5
$ws.Range("I23").Value = 225008.5
$ws.Range("K23").Value = 225008.5
$ws.Range("M23").Value = -224774.5
$ws.Range("H51").Value = 9199.6
$ws.Range("I51").Value = 8666.333000000001
$ws.Range("J51").Value = 9999.5
$ws.Range("K51").Value = 8666.333000000001
$ws.Range("L51").Value = 9999.5
$ws.Range("M51").Value = -8182.333000000001
$ws.Range("N51").Value = -10967.5
$ws.Range("H52").Value = 293.23077
$ws.Range("J52").Value = 293.23077
$ws.Range("L52").Value = 879.69231
$ws.Range("N52").Value = -1199.69231
$ws.Range("H86").Value = 2928165.8
$ws.Range("J86").Value = 4390690
$ws.Range("L86").Value = 4390690
$ws.Range("N86").Value = -4392936
$ws.Range("H89").Value = 2928165.8
$ws.Range("J89").Value = 4390690
$ws.Range("L89").Value = 21953450
$ws.Range("N89").Value = -21964682
$ws.Range("H107").Value = 129187
$ws.Range("I107").Value = 146928.14
$ws.Range("K107").Value = 146928.14
$ws.Range("M107").Value = -145008.14
$ws.Range("H131").Value = 4204.227
$ws.Range("I131").Value = 2535.7856
$ws.Range("J131").Value = 7124
$ws.Range("K131").Value = 7607.3568
$ws.Range("L131").Value = 21372
$ws.Range("M131").Value = -2567.3568
$ws.Range("N131").Value = -31452
$ws.Range("H132").Value = 3571.6875
$ws.Range("I132").Value = 3572.923
$ws.Range("K132").Value = 10718.769
$ws.Range("M132").Value = -8188.769

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7196.2666
$ws.Range("I45").Value = 2374.5
$ws.Range("K45").Value = 2374.5
$ws.Range("M45").Value = -1997.5
$ws.Range("H74").Value = 14181.458
$ws.Range("I74").Value = 17031.055
$ws.Range("J74").Value = 5632.6665
$ws.Range("K74").Value = 17031.055
$ws.Range("L74").Value = 5632.6665
$ws.Range("M74").Value = -16157.055
$ws.Range("N74").Value = -7380.6665
$ws.Range("H77").Value = 14181.458
$ws.Range("I77").Value = 17031.055
$ws.Range("J77").Value = 5632.6665
$ws.Range("K77").Value = 85155.27499999999
$ws.Range("L77").Value = 28163.3325
$ws.Range("M77").Value = -80787.27499999999
$ws.Range("N77").Value = -36899.3325
$ws.Range("H97").Value = 2781.4167
$ws.Range("I97").Value = 2781.4167
$ws.Range("K97").Value = 2781.4167
$ws.Range("M97").Value = -2285.4167

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 16999.5
$ws.Range("J88").Value = 16999.5
$ws.Range("L88").Value = 16999.5
$ws.Range("N88").Value = -17811.5
$ws.Range("H91").Value = 16999.5
$ws.Range("J91").Value = 16999.5
$ws.Range("L91").Value = 16999.5
$ws.Range("N91").Value = -19807.5
$ws.Range("H93").Value = 55000
$ws.Range("J93").Value = 55000
$ws.Range("L93").Value = 55000
$ws.Range("N93").Value = -58744
$ws.Range("H94").Value = 1650.375
$ws.Range("I94").Value = 1185.3077
$ws.Range("K94").Value = 1185.3077
$ws.Range("M94").Value = -734.3077000000001
$ws.Range("H105").Value = 1708.9231
$ws.Range("I105").Value = 990.25
$ws.Range("K105").Value = 990.25
$ws.Range("M105").Value = 756.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 446.2
$ws.Range("J7").Value = 349.6
$ws.Range("L7").Value = 349.6
$ws.Range("N7").Value = -575.6
$ws.Range("H41").Value = 31087.611
$ws.Range("J41").Value = 52414.285
$ws.Range("L41").Value = 52414.285
$ws.Range("N41").Value = -53270.285
$ws.Range("H50").Value = 28799.75
$ws.Range("J50").Value = 74950
$ws.Range("L50").Value = 74950
$ws.Range("N50").Value = -76200
$ws.Range("H51").Value = 65000
$ws.Range("I51").Value = 50000
$ws.Range("K51").Value = 50000
$ws.Range("M51").Value = -49264
$ws.Range("H60").Value = 93333.336
$ws.Range("J60").Value = 93333.336
$ws.Range("L60").Value = 93333.336
$ws.Range("N60").Value = -94355.336
$ws.Range("H61").Value = 65000
$ws.Range("I61").Value = 50000
$ws.Range("K61").Value = 50000
$ws.Range("M61").Value = -49652
$ws.Range("H62").Value = 3731.875
$ws.Range("J62").Value = 3471.75
$ws.Range("L62").Value = 3471.75
$ws.Range("N62").Value = -4719.75
$ws.Range("H65").Value = 3731.875
$ws.Range("J65").Value = 3471.75
$ws.Range("L65").Value = 17358.75
$ws.Range("N65").Value = -23598.75
$ws.Range("H68").Value = 120000
$ws.Range("J68").Value = 120000
$ws.Range("L68").Value = 120000
$ws.Range("N68").Value = -121498
$ws.Range("H71").Value = 120000
$ws.Range("J71").Value = 120000
$ws.Range("L71").Value = 360000
$ws.Range("N71").Value = -367488
$ws.Range("H74").Value = 120000
$ws.Range("J74").Value = 120000
$ws.Range("L74").Value = 120000
$ws.Range("N74").Value = -121748
$ws.Range("H77").Value = 120000
$ws.Range("J77").Value = 120000
$ws.Range("L77").Value = 360000
$ws.Range("N77").Value = -368736
$ws.Range("H122").Value = 2401.7
$ws.Range("I122").Value = 1496.6666
$ws.Range("K122").Value = 4489.9998
$ws.Range("M122").Value = -2039.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6892.552
$ws.Range("I102").Value = 8112.1
$ws.Range("K102").Value = 8112.1
$ws.Range("M102").Value = -6490.1
$ws.Range("H132").Value = 272966.5
$ws.Range("I132").Value = 335625.3
$ws.Range("K132").Value = 1006875.9
$ws.Range("M132").Value = -1004345.9

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3170.0435
$ws.Range("I46").Value = 2028.9286
$ws.Range("J46").Value = 4945.1113
$ws.Range("K46").Value = 2028.9286
$ws.Range("L46").Value = 4945.1113
$ws.Range("M46").Value = -1840.9286
$ws.Range("N46").Value = -5321.1113
$ws.Range("H61").Value = 4522.0586
$ws.Range("I61").Value = 3187.2222
$ws.Range("K61").Value = 3187.2222
$ws.Range("M61").Value = -2985.2222
$ws.Range("H82").Value = 2638.3635
$ws.Range("I82").Value = 2972.5
$ws.Range("K82").Value = 2972.5
$ws.Range("M82").Value = -2611.5
$ws.Range("H85").Value = 2638.3635
$ws.Range("I85").Value = 2972.5
$ws.Range("K85").Value = 2972.5
$ws.Range("M85").Value = -1724.5
$ws.Range("H113").Value = 4522.0586
$ws.Range("I113").Value = 3187.2222
$ws.Range("K113").Value = 3187.2222
$ws.Range("M113").Value = -1017.2222
$ws.Range("H122").Value = 492787.97
$ws.Range("I122").Value = 326371.78
$ws.Range("K122").Value = 979115.3400000001
$ws.Range("M122").Value = -976665.3400000001
$ws.Range("H132").Value = 3039.7917
$ws.Range("I132").Value = 1900.6364
$ws.Range("J132").Value = 4003.6924
$ws.Range("K132").Value = 5701.9092
$ws.Range("L132").Value = 12011.0772
$ws.Range("M132").Value = -3171.9092
$ws.Range("N132").Value = -17071.0772

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 252018.75
$ws.Range("I96").Value = 335358.34
$ws.Range("K96").Value = 335358.34
$ws.Range("M96").Value = -333985.34
$ws.Range("H113").Value = 1834.1875
$ws.Range("I113").Value = 1836.2
$ws.Range("K113").Value = 5508.6
$ws.Range("M113").Value = -3338.6
$ws.Range("H132").Value = 4046
$ws.Range("I132").Value = 4046
$ws.Range("K132").Value = 12138
$ws.Range("M132").Value = -9608
